{"js": "// Replace each two-digit-divided-by-one-digit expression with its new value.\n// Each \"old=\" string is unique in the document, so a plain (case-sensitive,\n// non-wildcard) search-and-replace is unambiguous and order-independent.\nconst replacements = [\n  [\"88\u00f79=\", \"69\u00f72=\"],\n  [\"98\u00f77=\", \"22\u00f77=\"],\n  [\"34\u00f79=\", \"51\u00f77=\"],\n  [\"25\u00f78=\", \"50\u00f77=\"],\n  [\"89\u00f72=\", \"85\u00f79=\"],\n  [\"47\u00f72=\", \"58\u00f75=\"],\n  [\"95\u00f76=\", \"52\u00f76=\"],\n  [\"75\u00f78=\", \"46\u00f76=\"],\n  [\"56\u00f73=\", \"78\u00f72=\"],\n  [\"98\u00f76=\", \"99\u00f74=\"],\n  [\"89\u00f75=\", \"34\u00f76=\"],\n  [\"42\u00f76=\", \"27\u00f74=\"],\n  [\"19\u00f79=\", \"76\u00f74=\"],\n  [\"72\u00f79=\", \"99\u00f79=\"],\n  [\"55\u00f75=\", \"41\u00f78=\"],\n  [\"75\u00f75=\", \"48\u00f73=\"],\n  [\"77\u00f78=\", \"28\u00f78=\"],\n  [\"17\u00f72=\", \"84\u00f73=\"],\n  [\"86\u00f76=\", \"24\u00f78=\"],\n  [\"57\u00f75=\", \"76\u00f72=\"],\n  [\"54\u00f73=\", \"66\u00f78=\"],\n  [\"81\u00f77=\", \"34\u00f72=\"],\n  [\"38\u00f72=\", \"81\u00f76=\"],\n  [\"72\u00f77=\", \"73\u00f79=\"],\n  [\"76\u00f74=\", \"75\u00f74=\"],\n];\n\n// Resolve every search FIRST (before any text is mutated) so that a\n// replacement's new text can never be accidentally re-matched by a\n// later rule searching the live document (e.g. \"19\u00f79=\" -> \"76\u00f74=\" must\n// not be re-caught by the unrelated \"76\u00f74=\" -> \"75\u00f74=\" rule).\nconst pending = [];\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  pending.push({ results, newText });\n}\nawait context.sync();\n\nfor (const { results, newText } of pending) {\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Replace each two-digit-divided-by-one-digit expression with its new value.\n# The document has a single table; the 5 data rows are table rows 1, 5, 9,\n# 13, 17 (the rows in between are blank spacer rows), each with 5 columns.\n# Addressing cells directly (rather than a live text search-and-replace)\n# sidesteps any \"domino\" risk where a later rule's search text could\n# accidentally match text that an earlier rule just inserted (e.g. rule\n# \"19\u00f79=\" -> \"76\u00f74=\" followed by unrelated rule \"76\u00f74=\" -> \"75\u00f74=\").\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$t.Cell(1,1).Range.Text  = \"69\u00f72=\"\n$t.Cell(1,2).Range.Text  = \"22\u00f77=\"\n$t.Cell(1,3).Range.Text  = \"51\u00f77=\"\n$t.Cell(1,4).Range.Text  = \"50\u00f77=\"\n$t.Cell(1,5).Range.Text  = \"85\u00f79=\"\n\n$t.Cell(5,1).Range.Text  = \"58\u00f75=\"\n$t.Cell(5,2).Range.Text  = \"52\u00f76=\"\n$t.Cell(5,3).Range.Text  = \"46\u00f76=\"\n$t.Cell(5,4).Range.Text  = \"78\u00f72=\"\n$t.Cell(5,5).Range.Text  = \"99\u00f74=\"\n\n$t.Cell(9,1).Range.Text  = \"34\u00f76=\"\n$t.Cell(9,2).Range.Text  = \"27\u00f74=\"\n$t.Cell(9,3).Range.Text  = \"76\u00f74=\"\n$t.Cell(9,4).Range.Text  = \"99\u00f79=\"\n$t.Cell(9,5).Range.Text  = \"41\u00f78=\"\n\n$t.Cell(13,1).Range.Text = \"48\u00f73=\"\n$t.Cell(13,2).Range.Text = \"28\u00f78=\"\n$t.Cell(13,3).Range.Text = \"84\u00f73=\"\n$t.Cell(13,4).Range.Text = \"24\u00f78=\"\n$t.Cell(13,5).Range.Text = \"76\u00f72=\"\n\n$t.Cell(17,1).Range.Text = \"66\u00f78=\"\n$t.Cell(17,2).Range.Text = \"34\u00f72=\"\n$t.Cell(17,3).Range.Text = \"81\u00f76=\"\n$t.Cell(17,4).Range.Text = \"73\u00f79=\"\n$t.Cell(17,5).Range.Text = \"75\u00f74=\"\n"}
